$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.349835276603699
$ws.Range("B1").Value = 2.689387798309326
$ws.Range("D1").Value = 1.549157023429871
$ws.Range("E1").Value = 0.9220970869064331
